$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds a new weekly price record for Mango at
# "Terminal Hortofrutícola Agro Chillán". It is inserted as the new
# row 175, pushing every following record down by one row
# (old row 175 -> 176, ..., old row 215 -> 216).
$ws.Rows(175).Insert()

$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 45258
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100108
$ws.Range("H175").Value = "Tropicales y subtropicales"
$ws.Range("I175").Value = 100108002
$ws.Range("J175").Value = "Mango"
$ws.Range("K175").Value = "Sin especificar"
$ws.Range("L175").Value = "Primera"
$ws.Range("M175").Value = 90
$ws.Range("N175").Value = 10000
$ws.Range("O175").Value = 11000
$ws.Range("P175").Value = 10556
$ws.Range("Q175").Value = "$/bandeja 4 kilos"
$ws.Range("R175").Value = "Brasil"
$ws.Range("S175").Value = 2639
$ws.Range("T175").Value = 4
